# Scottish Module Input files modification
# - "GeneralTaxRateMonthly" sheet/label renamed to "GeneralTaxRateWeekly"
# - "ProcessPayrollForMonthlyTax" sheet/label renamed to "ProcessPayrollForWeeklyTax"
# - "DO NOT TOUCH AUTOMATION EMP 105" test data renamed to "...EMP 107"
# - selections on each sheet updated to reflect where the author left the cursor

$wb = $excel.ActiveWorkbook

$wsFirst   = $wb.Worksheets.Item("first")
$wsMonthly = $wb.Worksheets.Item("GeneralTaxRateMonthly")
$wsProcess = $wb.Worksheets.Item("ProcessPayrollForMonthlyTax")
$wsReports = $wb.Worksheets.Item("TestReports")

# Update the employee test-data label (introduce the new shared string first so
# it lands in the same relative position as the source workbook).
$wsMonthly.Range("A2").Value = "DO NOT TOUCH AUTOMATION EMP 107"
$wsProcess.Range("B2").Value = "DO NOT TOUCH AUTOMATION EMP 107"
$wsReports.Range("B2").Value = "DO NOT TOUCH AUTOMATION EMP 107"

# Update the descriptive labels on the "first" sheet to match the renamed sheets.
$wsFirst.Range("A3").Value = "GeneralTaxRateWeekly"
$wsFirst.Range("A4").Value = "ProcessPayrollForWeeklyTax"

# Rename the sheets themselves (Weekly instead of Monthly).
$wsMonthly.Name = "GeneralTaxRateWeekly"
$wsProcess.Name = "ProcessPayrollForWeeklyTax"

# Re-point local references now that the sheets carry their new names.
$wsWeekly = $wb.Worksheets.Item("GeneralTaxRateWeekly")
$wsWeeklyProcess = $wb.Worksheets.Item("ProcessPayrollForWeeklyTax")

# Restore the per-sheet cursor/selection state left behind by the author.
$wsWeekly.Activate()
$wsWeekly.Range("C15").Select()

$wsWeeklyProcess.Activate()
$wsWeeklyProcess.Range("B2").Select()

$wsReports.Activate()
$wsReports.Range("B2").Select()

$wsFirst.Activate()
$wsFirst.Range("A3").Select()
